$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Osm"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 26.18311833333333
$ws.Range("H2").Value = 78.54935499999999
$ws.Range("I2").Value = 0.4541216129859197
$ws.Range("J2").Value = 0.4541216129859197
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 86.777428
$ws.Range("N2").Value = 173.554856
$ws.Range("O2").Value = 0.2676090626666408
$ws.Range("P2").Value = 0.2041884050300022
$ws.Range("Q2").Value = 2272.103665986313
$ws.Range("R2").Value = 13632.62199591788
$ws.Range("S2").Value = 0.121527059187825
$ws.Range("T2").Value = 0.09272636784524689

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Osm"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "Neutro"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 26.18311833333333
$ws.Range("H3").Value = 78.54935499999999
$ws.Range("I3").Value = 0.4541216129859197
$ws.Range("J3").Value = 0.4541216129859197
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 128.9086913333333
$ws.Range("N3").Value = 386.726074
$ws.Range("O3").Value = 0.3975357976419474
$ws.Range("P3").Value = 0.4549857149118007
$ws.Range("Q3").Value = 3375.231519375807
$ws.Range("R3").Value = 30377.08367438227
$ws.Range("S3").Value = 0.1805295976448053
$ws.Range("T3").Value = 0.2066188467412988

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Osm"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 26.18311833333333
$ws.Range("H4").Value = 78.54935499999999
$ws.Range("I4").Value = 0.4541216129859197
$ws.Range("J4").Value = 0.4541216129859197
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.52958833333333
$ws.Range("N4").Value = 37.588765
$ws.Range("O4").Value = 0.03863944192356349
$ws.Range("P4").Value = 0.04422342393230168
$ws.Range("Q4").Value = 328.0636939996194
$ws.Range("R4").Value = 2952.573245996575
$ws.Range("S4").Value = 0.01754700569120442
$ws.Range("T4").Value = 0.02008281260789696

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Osm"
$ws.Range("C5").Value = "Il6st"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 26.18311833333333
$ws.Range("H5").Value = 78.54935499999999
$ws.Range("I5").Value = 0.4541216129859197
$ws.Range("J5").Value = 0.4541216129859197
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.279662
$ws.Range("N5").Value = 51.83898599999999
$ws.Range("O5").Value = 0.05328798349515926
$ws.Range("P5").Value = 0.06098890064886812
$ws.Range("Q5").Value = 452.4354349060032
$ws.Range("R5").Value = 4071.918914154029
$ws.Range("S5").Value = 0.02419922501758879
$ws.Range("T5").Value = 0.027696377936902

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Osm"
$ws.Range("C6").Value = "Il6st"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 26.18311833333333
$ws.Range("H6").Value = 78.54935499999999
$ws.Range("I6").Value = 0.4541216129859197
$ws.Range("J6").Value = 0.4541216129859197
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 42.71737566666666
$ws.Range("N6").Value = 128.152127
$ws.Range("O6").Value = 0.1317342208129911
$ws.Range("P6").Value = 0.1507718021634167
$ws.Range("Q6").Value = 1118.474101969787
$ws.Range("R6").Value = 10066.26691772808
$ws.Range("S6").Value = 0.05982335684103883
$ws.Range("T6").Value = 0.06846873399124478

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Osm"
$ws.Range("C7").Value = "Il6st"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 26.18311833333333
$ws.Range("H7").Value = 78.54935499999999
$ws.Range("I7").Value = 0.4541216129859197
$ws.Range("J7").Value = 0.4541216129859197
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 36.0566465
$ws.Range("N7").Value = 72.113293
$ws.Range("O7").Value = 0.111193493459698
$ws.Range("P7").Value = 0.08484175331361067
$ws.Range("Q7").Value = 944.075442012669
$ws.Range("R7").Value = 5664.452652076015
$ws.Range("S7").Value = 0.05049536860345737
$ws.Range("T7").Value = 0.03852847386333038

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Osm"
$ws.Range("C8").Value = "Il6st"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 31.31438633333333
$ws.Range("H8").Value = 93.94315899999999
$ws.Range("I8").Value = 0.5431186404276995
$ws.Range("J8").Value = 0.5431186404276995
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 86.777428
$ws.Range("N8").Value = 173.554856
$ws.Range("O8").Value = 0.2676090626666408
$ws.Range("P8").Value = 0.2041884050300022
$ws.Range("Q8").Value = 2717.381905405017
$ws.Range("R8").Value = 16304.2914324301
$ws.Range("S8").Value = 0.145343470281637
$ws.Range("T8").Value = 0.1108985289309952

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Osm"
$ws.Range("C9").Value = "Il6st"
$ws.Range("D9").Value = "Neutro"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 31.31438633333333
$ws.Range("H9").Value = 93.94315899999999
$ws.Range("I9").Value = 0.5431186404276995
$ws.Range("J9").Value = 0.5431186404276995
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 128.9086913333333
$ws.Range("N9").Value = 386.726074
$ws.Range("O9").Value = 0.3975357976419474
$ws.Range("P9").Value = 0.4549857149118007
$ws.Range("Q9").Value = 4036.696562136418
$ws.Range("R9").Value = 36330.26905922776
$ws.Range("S9").Value = 0.2159091019366355
$ws.Range("T9").Value = 0.2471112228969221

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Osm"
$ws.Range("C10").Value = "Il6st"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 31.31438633333333
$ws.Range("H10").Value = 93.94315899999999
$ws.Range("I10").Value = 0.5431186404276995
$ws.Range("J10").Value = 0.5431186404276995
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.52958833333333
$ws.Range("N10").Value = 37.588765
$ws.Range("O10").Value = 0.03863944192356349
$ws.Range("P10").Value = 0.04422342393230168
$ws.Range("Q10").Value = 392.3563696676262
$ws.Range("R10").Value = 3531.207327008635
$ws.Range("S10").Value = 0.02098580116441085
$ws.Range("T10").Value = 0.02401856588116948

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Osm"
$ws.Range("C11").Value = "Il6st"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 31.31438633333333
$ws.Range("H11").Value = 93.94315899999999
$ws.Range("I11").Value = 0.5431186404276995
$ws.Range("J11").Value = 0.5431186404276995
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 17.279662
$ws.Range("N11").Value = 51.83898599999999
$ws.Range("O11").Value = 0.05328798349515926
$ws.Range("P11").Value = 0.06098890064886812
$ws.Range("Q11").Value = 541.1020115774193
$ws.Range("R11").Value = 4869.918104196773
$ws.Range("S11").Value = 0.02894169714702459
$ws.Range("T11").Value = 0.03312420880159329

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Osm"
$ws.Range("C12").Value = "Il6st"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 31.31438633333333
$ws.Range("H12").Value = 93.94315899999999
$ws.Range("I12").Value = 0.5431186404276995
$ws.Range("J12").Value = 0.5431186404276995
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 42.71737566666666
$ws.Range("N12").Value = 128.152127
$ws.Range("O12").Value = 0.1317342208129911
$ws.Range("P12").Value = 0.1507718021634167
$ws.Range("Q12").Value = 1337.668404772132
$ws.Range("R12").Value = 12039.01564294919
$ws.Range("S12").Value = 0.07154731090575409
$ws.Range("T12").Value = 0.08188697620582898

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Osm"
$ws.Range("C13").Value = "Il6st"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 31.31438633333333
$ws.Range("H13").Value = 93.94315899999999
$ws.Range("I13").Value = 0.5431186404276995
$ws.Range("J13").Value = 0.5431186404276995
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 36.0566465
$ws.Range("N13").Value = 72.113293
$ws.Range("O13").Value = 0.111193493459698
$ws.Range("P13").Value = 0.08484175331361067
$ws.Range("Q13").Value = 1129.091758385431
$ws.Range("R13").Value = 6774.550550312587
$ws.Range("S13").Value = 0.06039125899223748
$ws.Range("T13").Value = 0.0460791377111905

# Row 14
$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Osm"
$ws.Range("C14").Value = "Il6st"
$ws.Range("D14").Value = "FAPs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.1591176666666667
$ws.Range("H14").Value = 0.477353
$ws.Range("I14").Value = 0.002759746586380853
$ws.Range("J14").Value = 0.002759746586380853
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 86.777428
$ws.Range("N14").Value = 173.554856
$ws.Range("O14").Value = 0.2676090626666408
$ws.Range("P14").Value = 0.2041884050300022
$ws.Range("Q14").Value = 13.80782186269467
$ws.Range("R14").Value = 82.846931176168
$ws.Range("S14").Value = 0.0007385331971788417
$ws.Range("T14").Value = 0.0005635082537600995

# Row 15
$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Osm"
$ws.Range("C15").Value = "Il6st"
$ws.Range("D15").Value = "Neutro"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.1591176666666667
$ws.Range("H15").Value = 0.477353
$ws.Range("I15").Value = 0.002759746586380853
$ws.Range("J15").Value = 0.002759746586380853
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 128.9086913333333
$ws.Range("N15").Value = 386.726074
$ws.Range("O15").Value = 0.3975357976419474
$ws.Range("P15").Value = 0.4549857149118007
$ws.Range("Q15").Value = 20.51165017801356
$ws.Range("R15").Value = 184.604851602122
$ws.Range("S15").Value = 0.001097098060506554
$ws.Range("T15").Value = 0.001255645273579894

# Row 16
$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Osm"
$ws.Range("C16").Value = "Il6st"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.1591176666666667
$ws.Range("H16").Value = 0.477353
$ws.Range("I16").Value = 0.002759746586380853
$ws.Range("J16").Value = 0.002759746586380853
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 12.52958833333333
$ws.Range("N16").Value = 37.588765
$ws.Range("O16").Value = 0.03863944192356349
$ws.Range("P16").Value = 0.04422342393230168
$ws.Range("Q16").Value = 1.993678859893889
$ws.Range("R16").Value = 17.943109739045
$ws.Range("S16").Value = 0.0001066350679482155
$ws.Range("T16").Value = 0.0001220454432352429

# Row 17
$ws.Range("A17").Value = "Neutro"
$ws.Range("B17").Value = "Osm"
$ws.Range("C17").Value = "Il6st"
$ws.Range("D17").Value = "M2"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.1591176666666667
$ws.Range("H17").Value = 0.477353
$ws.Range("I17").Value = 0.002759746586380853
$ws.Range("J17").Value = 0.002759746586380853
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 17.279662
$ws.Range("N17").Value = 51.83898599999999
$ws.Range("O17").Value = 0.05328798349515926
$ws.Range("P17").Value = 0.06098890064886812
$ws.Range("Q17").Value = 2.749499498228667
$ws.Range("R17").Value = 24.745495484058
$ws.Range("S17").Value = 0.000147061330545885
$ws.Range("T17").Value = 0.0001683139103728348

# Row 18
$ws.Range("A18").Value = "Neutro"
$ws.Range("B18").Value = "Osm"
$ws.Range("C18").Value = "Il6st"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.1591176666666667
$ws.Range("H18").Value = 0.477353
$ws.Range("I18").Value = 0.002759746586380853
$ws.Range("J18").Value = 0.002759746586380853
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 42.71737566666666
$ws.Range("N18").Value = 128.152127
$ws.Range("O18").Value = 0.1317342208129911
$ws.Range("P18").Value = 0.1507718021634167
$ws.Range("Q18").Value = 6.797089142203444
$ws.Range("R18").Value = 61.17380227983099
$ws.Range("S18").Value = 0.0003635530661981937
$ws.Range("T18").Value = 0.0004160919663429786

# Row 19
$ws.Range("A19").Value = "Neutro"
$ws.Range("B19").Value = "Osm"
$ws.Range("C19").Value = "Il6st"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.1591176666666667
$ws.Range("H19").Value = 0.477353
$ws.Range("I19").Value = 0.002759746586380853
$ws.Range("J19").Value = 0.002759746586380853
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 36.0566465
$ws.Range("N19").Value = 72.113293
$ws.Range("O19").Value = 0.111193493459698
$ws.Range("P19").Value = 0.08484175331361067
$ws.Range("Q19").Value = 5.737249458904834
$ws.Range("R19").Value = 34.423496753429
$ws.Range("S19").Value = 0.0003068658640031633
$ws.Range("T19").Value = 0.0002341417390898035
